$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("ITR input data (2)")
$dst = $wb.Worksheets.Item("ITR input data")
$src.Rows(39).Copy() | Out-Null
$dst.Rows(33).PasteSpecial(-4104) | Out-Null
